# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets,
# which share identical data, per the commit's regenerated output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 2238
    $ws.Range("F3").Value = 1701
    $ws.Range("F4").Value = 332
    $ws.Range("F6").Value = 790
    $ws.Range("F8").Value = 5814
}
